# Insert a new data row for Halifax (YHZ) above the current Adelaide (ADL)
# row, shifting every following colo row down by one — matching the
# upstream "update generated data" regeneration that now lists YHZ first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 297:309 down to 298:310 and open up a blank row 297.
$ws.Rows(297).Insert()

# Populate the new row with the Halifax, Canada colo data.
$ws.Range("A297").Value = "YHZ"
$ws.Range("B297").Value = "Halifax, Canada"
$ws.Range("C297").Value = 44.64601
$ws.Range("D297").Value = -63.66844
$ws.Range("E297").Value = "CA"
$ws.Range("F297").Value = "North America"
$ws.Range("G297").Value = "Halifax"

# Match the formatting (bold, centered, bordered) used by every other
# "colo" cell in column A rather than the default Insert carry-down.
$ws.Range("A298").Copy()
$ws.Range("A297").PasteSpecial(-4122)
